{"js": "// Replace the full text of each of the 7 body paragraphs with the new\n// review content (title date, paper title, review paragraphs, link),\n// as described by the diff. Each paragraph holds a single run, so we\n// replace the whole paragraph's text while keeping paragraph formatting.\n\nconst newTexts = [\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 13.07.24: \u26a1\ufe0f\ud83d\ude80\",\n  \"SaySelf: Teaching LLMs to Express Confidence with Self-Reflective Rationales\",\n  \"\u05d1\u05d4\u05de\u05e9\u05da \u05dc\u05e1\u05e7\u05d9\u05e8\u05d4 \u05e9\u05dc \u05d0\u05ea\u05de\u05d5\u05dc, \u05de\u05d0\u05de\u05e8 \u05e7\u05dc\u05d9\u05dc \u05d9\u05d5\u05ea\u05e8 \u05e9\u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05dc\u05de\u05d3 \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05dc\u05e9\u05e2\u05e8\u05da \u05d0\u05d9 \u05d5\u05d5\u05d3\u05d0\u05d5\u05ea \u05d1\u05ea\u05e9\u05d5\u05d1\u05ea\u05dd. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d9\u05d8\u05d9\u05d1\u05d9\u05ea \u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05ea \u05de\u05e9\u05e0\u05d9 \u05e9\u05dc\u05d1\u05d9\u05dd \u05e2\u05d9\u05e7\u05e8\u05d9\u05d9\u05dd: \u05d9\u05e6\u05d9\u05e8\u05ea \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05dc\u05de\u05e9\u05d9\u05de\u05d4 \u05d6\u05d5 (\u05db\u05d9\u05de\u05d5\u05ea \u05d0\u05d9 \u05d5\u05d5\u05d3\u05d0\u05d5\u05ea) \u05d5\u05d8\u05d9\u05d5\u05d1 (fine-tuning) \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d4\u05d6\u05d4. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e9\u05e0\u05d9 \u05de\u05de\u05e9\u05d9\u05db\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05dd \u05e9\u05d9\u05d8\u05ea PPO \u05de\u05e2\u05d5\u05dc\u05dd \u05dc\u05de\u05d9\u05d3\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d7\u05d9\u05d6\u05d5\u05e7\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05e9\u05d9\u05e4\u05d5\u05e8 \u05e0\u05d5\u05e1\u05e3 \u05e9\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05d5.\",\n  \" \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05dc\u05d5\u05e7\u05d7\u05d9\u05dd \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05e9\u05d0\u05dc\u05d5\u05ea \u05d5\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05d4\u05e0\u05e7\u05e8\u05d0 HotpotQA \u05d5\u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea \u05de\u05de\u05e0\u05d5 \u05dc\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05d5\u05de\u05d1\u05e7\u05e9\u05d9\u05dd \u05de\u05de\u05e0\u05d5 \u05dc\u05ea\u05ea \u05ea\u05e9\u05d5\u05d1\u05d4 \u05de\u05dc\u05d5\u05d5\u05d4 \u05d1-reasoning. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05e7\u05dc\u05e1\u05d8\u05e8\u05d9\u05dd \u05d0\u05ea \u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05d4\u05de\u05d5\u05d3\u05dc (\u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4-reasoning) \u05dc\u05e7\u05dc\u05e1\u05d8\u05e8\u05d9\u05dd \u05dc\u05e4\u05d9 \u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc\u05d4\u05dd \u05d5\u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d9\u05d7\u05e1 \u05e9\u05dc \u05d2\u05d5\u05d3\u05dc \u05d4\u05e7\u05dc\u05e1\u05d8\u05e8 \u05d4\u05de\u05db\u05d9\u05dc \u05d0\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e0\u05db\u05d5\u05e0\u05d4 (\u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8) \u05d9\u05d7\u05e1\u05d9\u05ea \u05dc\u05db\u05dc \u05d4\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea. \u05d6\u05d4 \u05de\u05d3\u05d3 \u05d0\u05d9 \u05d4\u05d5\u05d5\u05d3\u05d0\u05d5\u05ea \u05e9\u05dc\u05e0\u05d5 \u05e9\u05e2\u05dc\u05d9\u05d5 \u05e0\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05d1\u05d4\u05de\u05e9\u05da.\",\n  \"\u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05e4\u05dc\u05d8\u05e8\u05d9\u05dd \u05d0\u05ea \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea \u05d5\u05d1\u05e1\u05d5\u05e3 \u05de\u05d1\u05e7\u05e9\u05d9\u05dd \u05de-gpt4 \u05dc\u05ea\u05ea \u05d4\u05e1\u05d1\u05e8\u05d9\u05dd \u05dc\u05de\u05d4 \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d9\u05d4 \u05e2\u05e9\u05d5\u05d9 \u05dc\u05ea\u05ea \u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05dc\u05d0 \u05e0\u05db\u05d5\u05e0\u05d5\u05ea \u05dc\u05e9\u05d0\u05dc\u05d4 (\u05db\u05dc\u05d5\u05de\u05e8 \u05f4\u05d4\u05e1\u05d9\u05d1\u05d4\u05f4 \u05dc\u05d0\u05d9 \u05d5\u05d5\u05d3\u05d0\u05d5\u05ea). \u05d1\u05e9\u05dc\u05d1 \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05de\u05d8\u05d9\u05d9\u05d1\u05d9\u05dd (\u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05e0\u05ea\u05d5\u05df) \u05e7\u05d5\u05d3\u05d9\u05dd \u05db\u05dc \u05dc\u05ea\u05ea \u05ea\u05e9\u05d5\u05d1\u05d4 \u05e0\u05db\u05d5\u05df, \u05dc\u05d3\u05d9\u05d9\u05e7 \u05d1\u05de\u05de\u05d3 \u05e9\u05dc \u05d0\u05d9 \u05d4\u05d5\u05d5\u05d3\u05d0\u05d5\u05ea \u05d5\u05d1\u05e0\u05d5\u05e1\u05e3 \u05dc\u05ea\u05ea reasoning \u05e0\u05db\u05d5\u05df \u05dc\u05e0\u05d5\u05db\u05d7\u05d5\u05ea \u05e9\u05dc \u05d0\u05d9 \u05d4\u05d5\u05d5\u05d3\u05d0\u05d5\u05ea. \u05db\u05dc \u05d0\u05dc\u05d4 \u05e0\u05de\u05e6\u05d0\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05d5 \u05de\u05e4\u05d5\u05e8\u05e9 \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05dc\u05d5\u05e1.\",\n  \"\u05d1\u05e9\u05dc\u05d1 \u05d4\u05e9\u05e0\u05d9 \u05de\u05de\u05e9\u05d9\u05db\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05d1\u05e9\u05d9\u05d8\u05d4 PPO \u05db\u05d3\u05d9 \u05dc\u05de\u05d6\u05e2\u05e8 (\u05d0\u05d5 \u05dc\u05de\u05e7\u05e1\u05dd \u05d0\u05d5\u05ea\u05d4 \u05e2\u05dd \u05de\u05d9\u05e0\u05d5\u05e1) \u05d0\u05ea \u05d4\u05d4\u05e4\u05e8\u05e9 \u05d1\u05d9\u05df \u05e0\u05db\u05d5\u05e0\u05d5\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 (0 \u05d0\u05d5 1) \u05d5\u05e8\u05de\u05ea \u05d4-confidence \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05dc\u05d2\u05d1\u05d9\u05d4. \u05db\u05de\u05d5 \u05d1\u05db\u05dc \u05e9\u05d9\u05d8\u05ea PPO \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e0\u05d5\u05e6\u05e8\u05d5\u05ea \u05f4on the fly\\\" \u05d0\u05d7\u05e8\u05d9 \u05db\u05dc \u05e2\u05d3\u05db\u05d5\u05df \u05e9\u05dc \u05de\u05e9\u05e7\u05dc\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc.\",\n  \"https://arxiv.org/abs/2405.20974\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newTexts.length) {\n  throw new Error(\n    `Unexpected paragraph count: expected ${newTexts.length}, found ${paragraphs.items.length}`\n  );\n}\n\nfor (let i = 0; i < newTexts.length; i++) {\n  paragraphs.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the full text of each of the 7 body paragraphs with the new\n# review content (title date, paper title, review paragraphs, link), as\n# described by the diff. Each paragraph's Range.Text is reassigned in\n# place (Word keeps the trailing paragraph mark automatically), so\n# paragraph formatting/count is preserved.\n\n$d = $word.ActiveDocument\n\n$newTexts = @(\n    '\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 13.07.24: \u26a1\ufe0f\ud83d\ude80',\n    'SaySelf: Teaching LLMs to Express Confidence with Self-Reflective Rationales',\n    '\u05d1\u05d4\u05de\u05e9\u05da \u05dc\u05e1\u05e7\u05d9\u05e8\u05d4 \u05e9\u05dc \u05d0\u05ea\u05de\u05d5\u05dc, \u05de\u05d0\u05de\u05e8 \u05e7\u05dc\u05d9\u05dc \u05d9\u05d5\u05ea\u05e8 \u05e9\u05de\u05e6\u05d9\u05e2 \u05e9\u05d9\u05d8\u05d4 \u05dc\u05dc\u05de\u05d3 \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05dc\u05e9\u05e2\u05e8\u05da \u05d0\u05d9 \u05d5\u05d5\u05d3\u05d0\u05d5\u05ea \u05d1\u05ea\u05e9\u05d5\u05d1\u05ea\u05dd. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d9\u05d8\u05d9\u05d1\u05d9\u05ea \u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05ea \u05de\u05e9\u05e0\u05d9 \u05e9\u05dc\u05d1\u05d9\u05dd \u05e2\u05d9\u05e7\u05e8\u05d9\u05d9\u05dd: \u05d9\u05e6\u05d9\u05e8\u05ea \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05dc\u05de\u05e9\u05d9\u05de\u05d4 \u05d6\u05d5 (\u05db\u05d9\u05de\u05d5\u05ea \u05d0\u05d9 \u05d5\u05d5\u05d3\u05d0\u05d5\u05ea) \u05d5\u05d8\u05d9\u05d5\u05d1 (fine-tuning) \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05dc \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05d4\u05d6\u05d4. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e9\u05e0\u05d9 \u05de\u05de\u05e9\u05d9\u05db\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05e2\u05dd \u05e9\u05d9\u05d8\u05ea PPO \u05de\u05e2\u05d5\u05dc\u05dd \u05dc\u05de\u05d9\u05d3\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d7\u05d9\u05d6\u05d5\u05e7\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05e9\u05d9\u05e4\u05d5\u05e8 \u05e0\u05d5\u05e1\u05e3 \u05e9\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05d5.',\n    ' \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05dc\u05d5\u05e7\u05d7\u05d9\u05dd \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05e9\u05d0\u05dc\u05d5\u05ea \u05d5\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05d4\u05e0\u05e7\u05e8\u05d0 HotpotQA \u05d5\u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea \u05de\u05de\u05e0\u05d5 \u05dc\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05d5\u05de\u05d1\u05e7\u05e9\u05d9\u05dd \u05de\u05de\u05e0\u05d5 \u05dc\u05ea\u05ea \u05ea\u05e9\u05d5\u05d1\u05d4 \u05de\u05dc\u05d5\u05d5\u05d4 \u05d1-reasoning. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05e7\u05dc\u05e1\u05d8\u05e8\u05d9\u05dd \u05d0\u05ea \u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05d4\u05de\u05d5\u05d3\u05dc (\u05d9\u05d7\u05d3 \u05e2\u05dd \u05d4-reasoning) \u05dc\u05e7\u05dc\u05e1\u05d8\u05e8\u05d9\u05dd \u05dc\u05e4\u05d9 \u05d4\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2 \u05e9\u05dc\u05d4\u05dd \u05d5\u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d9\u05d7\u05e1 \u05e9\u05dc \u05d2\u05d5\u05d3\u05dc \u05d4\u05e7\u05dc\u05e1\u05d8\u05e8 \u05d4\u05de\u05db\u05d9\u05dc \u05d0\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e0\u05db\u05d5\u05e0\u05d4 (\u05de\u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8) \u05d9\u05d7\u05e1\u05d9\u05ea \u05dc\u05db\u05dc \u05d4\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea. \u05d6\u05d4 \u05de\u05d3\u05d3 \u05d0\u05d9 \u05d4\u05d5\u05d5\u05d3\u05d0\u05d5\u05ea \u05e9\u05dc\u05e0\u05d5 \u05e9\u05e2\u05dc\u05d9\u05d5 \u05e0\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05d1\u05d4\u05de\u05e9\u05da.',\n    '\u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05e4\u05dc\u05d8\u05e8\u05d9\u05dd \u05d0\u05ea \u05d4\u05e9\u05d0\u05dc\u05d5\u05ea \u05d5\u05d1\u05e1\u05d5\u05e3 \u05de\u05d1\u05e7\u05e9\u05d9\u05dd \u05de-gpt4 \u05dc\u05ea\u05ea \u05d4\u05e1\u05d1\u05e8\u05d9\u05dd \u05dc\u05de\u05d4 \u05d4\u05de\u05d5\u05d3\u05dc \u05d4\u05d9\u05d4 \u05e2\u05e9\u05d5\u05d9 \u05dc\u05ea\u05ea \u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05dc\u05d0 \u05e0\u05db\u05d5\u05e0\u05d5\u05ea \u05dc\u05e9\u05d0\u05dc\u05d4 (\u05db\u05dc\u05d5\u05de\u05e8 \u05f4\u05d4\u05e1\u05d9\u05d1\u05d4\u05f4 \u05dc\u05d0\u05d9 \u05d5\u05d5\u05d3\u05d0\u05d5\u05ea). \u05d1\u05e9\u05dc\u05d1 \u05d4\u05d0\u05d7\u05e8\u05d5\u05df \u05de\u05d8\u05d9\u05d9\u05d1\u05d9\u05dd (\u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05e0\u05ea\u05d5\u05df) \u05e7\u05d5\u05d3\u05d9\u05dd \u05db\u05dc \u05dc\u05ea\u05ea \u05ea\u05e9\u05d5\u05d1\u05d4 \u05e0\u05db\u05d5\u05df, \u05dc\u05d3\u05d9\u05d9\u05e7 \u05d1\u05de\u05de\u05d3 \u05e9\u05dc \u05d0\u05d9 \u05d4\u05d5\u05d5\u05d3\u05d0\u05d5\u05ea \u05d5\u05d1\u05e0\u05d5\u05e1\u05e3 \u05dc\u05ea\u05ea reasoning \u05e0\u05db\u05d5\u05df \u05dc\u05e0\u05d5\u05db\u05d7\u05d5\u05ea \u05e9\u05dc \u05d0\u05d9 \u05d4\u05d5\u05d5\u05d3\u05d0\u05d5\u05ea. \u05db\u05dc \u05d0\u05dc\u05d4 \u05e0\u05de\u05e6\u05d0\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05d5 \u05de\u05e4\u05d5\u05e8\u05e9 \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05dc\u05d5\u05e1.',\n    '\u05d1\u05e9\u05dc\u05d1 \u05d4\u05e9\u05e0\u05d9 \u05de\u05de\u05e9\u05d9\u05db\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05d1\u05e9\u05d9\u05d8\u05d4 PPO \u05db\u05d3\u05d9 \u05dc\u05de\u05d6\u05e2\u05e8 (\u05d0\u05d5 \u05dc\u05de\u05e7\u05e1\u05dd \u05d0\u05d5\u05ea\u05d4 \u05e2\u05dd \u05de\u05d9\u05e0\u05d5\u05e1) \u05d0\u05ea \u05d4\u05d4\u05e4\u05e8\u05e9 \u05d1\u05d9\u05df \u05e0\u05db\u05d5\u05e0\u05d5\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 (0 \u05d0\u05d5 1) \u05d5\u05e8\u05de\u05ea \u05d4-confidence \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc \u05dc\u05d2\u05d1\u05d9\u05d4. \u05db\u05de\u05d5 \u05d1\u05db\u05dc \u05e9\u05d9\u05d8\u05ea PPO \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e0\u05d5\u05e6\u05e8\u05d5\u05ea \u05f4on the fly\" \u05d0\u05d7\u05e8\u05d9 \u05db\u05dc \u05e2\u05d3\u05db\u05d5\u05df \u05e9\u05dc \u05de\u05e9\u05e7\u05dc\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc.',\n    'https://arxiv.org/abs/2405.20974'\n)\n\nif ($d.Paragraphs.Count -ne $newTexts.Count) {\n    throw \"Unexpected paragraph count: expected $($newTexts.Count), found $($d.Paragraphs.Count)\"\n}\n\nfor ($i = 0; $i -lt $newTexts.Count; $i++) {\n    $d.Paragraphs($i + 1).Range.Text = $newTexts[$i]\n}\n"}
